$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -11.4683
$ws.Range("B3").Value = 5.98309999999999
$ws.Range("C5").Value = -13.56579999999999
$ws.Range("B14").Value = 9.121700000000001
$ws.Range("B16").Value = 9.324600000000002
$ws.Range("C16").Value = -11.64620000000001
$ws.Range("B21").Value = 5.651499999999996
$ws.Range("B23").Value = 5.545
$ws.Range("B25").Value = 5.828299999999995
